# Apply "Add data for 2022-05-05" update to the carjacking-by-neighborhood-by-month workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet to reflect the new "through" date.
$ws.Name = "Through 2022-04-27"

# Update the header label for the current-month column (column B).
$ws.Range("B1").Value = "April 2022 (through April 27)"

# Update/insert individual neighborhood x April-year cell counts.
$ws.Range("N2").Value = 7     # Austin, April 2019
$ws.Range("Z2").Value = 2     # Austin, April 2016
$ws.Range("V3").Value = 4     # Englewood, April 2017
$ws.Range("AD4").Value = 1    # North Lawndale, April 2015
$ws.Range("AD5").Value = 3    # Garfield Park, April 2015
$ws.Range("J6").Value = 2     # Humboldt Park, April 2020
$ws.Range("Z13").Value = 1    # Wicker Park, April 2016
$ws.Range("F25").Value = 4    # South Shore, April 2021
$ws.Range("F39").Value = 3    # Little Village, April 2021
$ws.Range("F44").Value = 1    # Brighton Park, April 2021
$ws.Range("AD50").Value = 1   # Grand Crossing, April 2015
$ws.Range("J89").Value = 2    # South Chicago, April 2020
$ws.Range("B92").Value = 1    # United Center, April 2022 (through April 27)
